$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-Text "Unraveling the Enigma of Dark Matter" "Unveiling the Intricate World of Chemistry: Exploring the Essence of Matter"

# ---------------------------------------------------------------------------
# 2. Author name
# ---------------------------------------------------------------------------
Replace-Text "Damien Brooks" "Dr. Emily Sanders"

# ---------------------------------------------------------------------------
# 3. Email
# ---------------------------------------------------------------------------
Replace-Text "damien" "esanders@schoolmail"
Replace-Text "brooks@darkmatter." ""
Replace-Text "net" "edu"

# ---------------------------------------------------------------------------
# 4. Body paragraph 1 (intro)
# ---------------------------------------------------------------------------
Replace-Text "The cosmos, vast and enigmatic, holds secrets that intrigue the most astute minds" "Imagine a world where everything is made up of tiny particles called atoms and molecules, interacting in a complex dance to create the universe around us"

Replace-Text " Among these mysteries, dark matter looms as one of the most perplexing enigmas" " Chemistry is the study of these interactions, a science that delves into the structure, properties, and behavior of matter"

Replace-Text " This elusive substance, believed to constitute approximately 27% of the universe, exerts a gravitational pull, shaping galaxies and influencing the motion of celestial objects. Yet, despite its profound influence, dark matter remains shrouded in mystery." ""

Replace-Text " As we embark on this intellectual journey, we will delve into the intricate nature of dark matter and contemplate its implications for our understanding of the universe" " As we embark on this journey through the realm of chemistry, we will unravel the secrets of the elements, discover the wonders of chemical reactions, and uncover the pervasive influence of chemistry in our everyday lives"

Replace-Text "In 1933, Swiss astrophysicist Fritz Zwicky conducted a meticulous analysis of the Coma Cluster, a dense gathering of galaxies" "In the tapestry of chemistry, we will encounter elements, the fundamental building blocks of matter, each possessing unique properties that contribute to the diversity of the world around us"

Replace-Text " Through intricate calculations, he discovered a discrepancy between the observed velocity of galaxies and the mass estimated from visible matter alone" " We will delve into the periodic table, a treasure map of elements, organized by their atomic number and revealing patterns that govern their reactivity and behavior"

Replace-Text " This discrepancy hinted at the presence of an invisible mass, later termed dark matter" " Through chemical reactions, we will witness the transformation of substances, as atoms rearrange themselves to form new substances with different properties"

Replace-Text " Subsequent observations reinforced Zwicky's findings, confirming the substantial contribution of dark matter to the universe's overall mass." ""

Replace-Text " Despite these tantalizing clues, dark matter's elusive nature has continued to evade direct detection" " We will explore the concepts of energy transfer and chemical bonding, understanding how energy drives chemical reactions and how atoms unite to form molecules"

Replace-Text "The search for dark matter has become a scientific endeavor of paramount importance" "Chemistry extends far beyond the confines of the laboratory"

Replace-Text " Scientists have employed various techniques to unravel its enigmatic properties" " It plays a crucial role in diverse fields, influencing everything from medicine and agriculture to energy and materials science"

Replace-Text " Underground laboratories shield sensitive detectors from cosmic radiation, seeking to capture the faint signals of dark matter particles" " In the realm of medicine, chemistry enables the development of drugs to combat diseases, while in agriculture, it aids in the creation of fertilizers and pesticides to enhance crop yields"

Replace-Text " Space telescopes scan the universe, observing gravitational lensing effects caused by the presence of dark matter. Particle accelerators, such as the Large Hadron Collider, collide particles at high energies, hoping to produce dark matter particles that can be detected." ""

Replace-Text " Though these efforts have yielded valuable insights, the true nature of dark matter remains an enigma, challenging our understanding of the fundamental laws that govern the universe" " Chemistry also empowers us to harness energy from various sources, from fossil fuels to renewable resources, and to create advanced materials with tailored properties"

# ---------------------------------------------------------------------------
# 5. Summary paragraph
# ---------------------------------------------------------------------------
Replace-Text "Dark matter, an enigmatic substance constituting a significant portion of the universe's mass, remains an enduring mystery" "Chemistry is a captivating science that delves into the intricacies of matter, revealing the fundamental principles governing the universe around us"

Replace-Text " Despite its profound gravitational influence, its elusive nature has thwarted attempts at direct detection" " By exploring the world of elements, compounds, and reactions, we unveil the secrets of chemical transformations and gain insights into the pervasive influence of chemistry in our lives"

Replace-Text " The quest to unravel the secrets of dark matter continues with innovative experimental techniques and theoretical models. Its discovery " ""

Replace-Text "promises to revolutionize our understanding of cosmology and deepen our knowledge of the fundamental forces that shape the universe" " From the study of the periodic table to the examination of energy transfer and chemical bonding, chemistry provides a lens through which we can comprehend the complex interactions shaping our world"

# ---------------------------------------------------------------------------
# 6. Trailing empty paragraph
# ---------------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
